$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_QueryLog_External")

# --- Update "Effort (PDs) Incl. of Testing" (N) and "Comments - On PD(s)" (O) columns ---
# These cells previously held a shared formula (=M*1.5) but are now overwritten with
# literal text values captured from the updated change log.

$ws.Range("N11").Value = "24"
$ws.Range("O11").Value = "Assumptions : RP will expose 2 apis`n1. to provide applicant info.`n2. receive updated info and incorporate it."

$ws.Range("N12").Value = "36"
$ws.Range("O12").Value = "Estimation may change after understanding overall scope of the change."

$ws.Range("N13").Value = "52"
$ws.Range("O13").Value = "Need more clarification on the requirement. Estimation may change after clarification."

$ws.Range("N14").Value = "30"

$ws.Range("N15").Value = "28"
$ws.Range("O15").Value = "Estimation may change after understanding overall scope of the change."

$ws.Range("N20").Value = "20"
$ws.Range("O20").Value = "Since the requirement is not detailed the effort may change."

$ws.Range("N40").Value = "45"

$ws.Range("N41").Value = "55"
$ws.Range("O41").Value = "Need more clarification on the requirement. Estimation may change after clarification."

$ws.Range("N42").Value = "12"
$ws.Range("O42").Value = "Change algorithm from lavenstine distance to phonetic and soundex match"

# --- Re-apply the AutoFilter over the full data range, filtered to Module = "Registration Processor" ---

$ws.AutoFilterMode = $false
$ws.Range("A2:H53").AutoFilter()
$ws.Range("A2:H53").AutoFilter(4, @("Registration Processor"))

# Keep the workbook-level _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "MOSIP_QueryLog_External!_FilterDatabase") {
        $n.RefersTo = "=MOSIP_QueryLog_External!`$A`$2:`$H`$53"
    }
}
